# Applies the "cryptos list" refresh described in the commit message:
# updates Price (D) / Volume(1h) (E) figures, and reorders the
# InternetComputer(DFINITY) / RenderToken rows (29 <-> 30).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value2 = "65.609.12"
$ws.Range("E2").Value2 = "  -0.86%  "
# Row 3
$ws.Range("D3").Value2 = "3.437.60"
$ws.Range("E3").Value2 = "  -3.04%  "
# Row 4
$ws.Range("E4").Value2 = "  +0.00%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = "590.25"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value2 = "  -2.09%  "
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value2 = "137.98"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value2 = "  -5.94%  "
# Row 7
$ws.Range("D7").Value2 = "3.439.30"
$ws.Range("E7").Value2 = "  -3.00%  "
# Row 8
$ws.Range("E8").Value2 = "  +0.04%  "
# Row 9
$ws.Range("E9").Value2 = "  -0.10%  "
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value2 = "7.34"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value2 = "  -6.00%  "
# Row 11
$ws.Range("E11").Value2 = "  -8.62%  "
# Row 13
$ws.Range("D13").Value2 = "4.019.26"
$ws.Range("E13").Value2 = "  -2.99%  "
# Row 14
$ws.Range("E14").Value2 = "  -10.20%  "
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value2 = "26.49"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value2 = "  -8.85%  "
# Row 16
$ws.Range("D16").Value2 = "3.434.66"
$ws.Range("E16").Value2 = "  -2.83%  "
# Row 17
$ws.Range("D17").Value2 = "65.532.88"
$ws.Range("E17").Value2 = "  -0.90%  "
# Row 18
$ws.Range("E18").Value2 = "  -1.67%  "
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value2 = "9.87"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value2 = "  -10.63%  "
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value2 = "5.89"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value2 = "  -5.67%  "
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value2 = "13.70"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value2 = "  -6.59%  "
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value2 = "392.76"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value2 = "  -6.44%  "
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value2 = "0.554"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value2 = "  -8.07%  "
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value2 = "73.38"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value2 = "  -5.76%  "
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value2 = "0.999"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value2 = "  -0.09%  "
# Row 26
$ws.Range("D26").Value2 = "3.573.55"
$ws.Range("E26").Value2 = "  -2.94%  "
# Row 27
$ws.Range("E27").Value2 = "  -8.09%  "
# Row 28
$ws.Range("E28").Value2 = "  +0.22%  "
# Row 29
$ws.Range("B29").Value2 = "RenderToken"
$ws.Range("C29").Value2 = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value2 = "7.20"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value2 = "  -8.06%  "
# Row 30
$ws.Range("B30").Value2 = "InternetComputer(DFINITY)"
$ws.Range("C30").Value2 = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value2 = "8.27"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value2 = "  -9.57%  "
# Row 31
$ws.Range("E31").Value2 = "  -9.36%  "
# Row 32
$ws.Range("D32").Value2 = "3.442.16"
$ws.Range("E32").Value2 = "  -2.74%  "
# Row 33
$ws.Range("E33").Value2 = "  -0.02%  "
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value2 = "0.145"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value2 = "  -6.88%  "
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value2 = "23.02"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value2 = "  -5.97%  "
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value2 = "172.37"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value2 = "  -1.04%  "
# Row 37
$ws.Range("E37").Value2 = "  -9.31%  "
# Row 38
$ws.Range("E38").Value2 = "  -9.41%  "
# Row 39
$ws.Range("E39").Value2 = "  -7.81%  "
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value2 = "4.81"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value2 = "  -8.99%  "
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value2 = "0.0763"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value2 = "  -7.71%  "
# Row 42
$ws.Range("E42").Value2 = "  -4.28%  "
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value2 = "43.77"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value2 = "  -4.10%  "
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value2 = "1.00"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value2 = "  +0.09%  "
# Row 45
$ws.Range("E45").Value2 = "  -13.47%  "
# Row 46
$ws.Range("E46").Value2 = "  -10.15%  "
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value2 = "1.12"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value2 = "  +1.77%  "
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value2 = "22.75"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value2 = "  +0.20%  "
# Row 49
$ws.Range("E49").Value2 = "  -8.05%  "
# Row 50
$ws.Range("E50").Value2 = "  -12.63%  "
# Row 51
$ws.Range("D51").Value2 = "2.211.80"
$ws.Range("E51").Value2 = "  -6.97%  "
